$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh: update Price (D) and Volume(1h) (E) columns for each row,
# and for rows whose coin was re-ranked, the Coin (B) / Link (C) text as well.
#
# D4 is an untouched data cell with the default (General) number format/style;
# we snapshot its .Style so that after forcing text entry (leading "'") on
# numeric-looking price strings (to stop Excel re-parsing "543.85" -> 543.85000000000002,
# "1.00" -> 1, etc.), we can restore the cell to that same unstyled look.
$normalStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = '59.535.19'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '2.530.33'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = "'543.85"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = "'145.61"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").Value = '2.557.75'
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = "'5.59"
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").Value = "'0.360"
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '2.973.17'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = "'23.61"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = '  -4.36%  '
$ws.Range("D16").Value = '59.445.44'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '2.539.87'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").Value = "'11.24"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").Value = "'4.29"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").Value = "'326.96"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").Value = "'5.94"
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").Value = "'62.41"
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").Value = "'0.437"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = '  -2.76%  '
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").Value = "'0.994"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").Value = "'8.02"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("D29").Value = '0.0₃0794'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = "'6.79"
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.82"
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("D32").Value = "'1.21"
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = '  -10.21%  '
$ws.Range("D33").Value = "'1.49"
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").Value = "'161.08"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").Value = "'18.77"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").Value = "'4.41"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = '  -3.17%  '
$ws.Range("D38").Value = "'1.62"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  -7.93%  '
$ws.Range("D39").Value = "'37.16"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("D40").Value = "'5.63"
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = '  -7.79%  '
$ws.Range("D41").Value = "'0.843"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("D42").Value = "'297.40"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  -6.19%  '
$ws.Range("D43").Value = "'3.71"
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'0.609"
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = "'0.993"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("D46").Value = "'10.81"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'18.98"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = "'0.0939"
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").Value = "'123.80"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = '  -2.62%  '
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").Value = "'0.0516"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0228"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = '  -2.23%  '
